$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- 1. Remove the "Año de la recuperación..." run and its paragraph's
#        centered alignment from the default header (header2.xml). ---
$hdr = $sec.Headers.Item(1)  # wdHeaderFooterPrimary -> default header
$rng = $hdr.Range
$targetText = [string]::Concat([char]8220, "Año de la recuperación y consolidación de la economía peruana", [char]8221)

$found = $rng.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs.Item(1)
    $para.Alignment = 0   # wdAlignParagraphLeft -> drops <w:jc w:val="center"/>
}

$rng2 = $hdr.Range
$rng2.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- 2. Add even/default/first footers, each a single empty paragraph
#        styled "Piedepgina" (mirrors the existing footer style). ---
$styleCountBefore = $d.Styles.Count

for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    $ftr.Range.Text = ""
    $p = $ftr.Range.Paragraphs.Item(1)
    $p.Range.Style = "Piedepgina"
}

# Creating the footer parts mints fresh built-in "Header"/"Footer" style
# pairs (even though equivalent localized styles already exist). Remove
# just those newly minted trailing styles so styles.xml stays untouched.
for ($i = $d.Styles.Count; $i -gt $styleCountBefore; $i--) {
    $d.Styles.Item($i).Delete()
}
